# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamp values that get refreshed each time
# the handback status report is (re)generated.
#
# These cells hold plain text timestamps (format "yyyy-mm-dd HH:mm:ss") stored
# as shared strings, so we assign them as literal strings (not as Date values)
# to avoid Excel re-typing them as numeric dates.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 7f45ec9f-...md (row 4)
$overview.Range("G4").Value = "2016-08-13 12:57:40"

# zh-cn sheet: handoff / handback datetimes for 7f45ec9f-...xlf (row 4)
$zhcn.Range("H4").Value = "2016-08-13 12:57:32"
$zhcn.Range("K4").Value = "2016-08-13 12:58:03"

# de-de sheet: handback datetime for 7f45ec9f-...xlf (row 4)
# (its "Correspond Handoff Datetime" H4 shares the same value as Overview!G4
# and is already updated above via the shared string)
$dede.Range("K4").Value = "2016-08-13 12:58:14"
